$d = $word.ActiveDocument

# --- 1. Update the visible sentence text (formatting of the existing runs is
#        preserved automatically because Find/Replace edits text in place). ---
$d.Content.Find.Execute("This sentence is inserted on 13", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Here comes in a sentence inserted on 14", 2)

$d.Content.Find.Execute(" October, 2013.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " October, 2018", 2)

# --- 2. Locate the paragraph that now ends in " October, 2018" followed by the
#        _GoBack bookmark, and append a new run containing "." right after the
#        bookmark (matching a run that was re-typed after the bookmark). ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*October, 2018*") {
        $target = $cand
    }
}

$paraEnd = $target.Range.End - 1   # right after the _GoBack bookmark, before the pilcrow

# Borrow the run formatting already used in this paragraph (snapToGrid/kern) by
# duplicating one of its characters as FormattedText, pasting it after the
# bookmark, then swapping its text for the period - this keeps the new run's
# rPr identical to its neighbours without disturbing any other run.
$donor = $d.Range($target.Range.Start, $target.Range.Start + 1)
$formatted = $donor.FormattedText

$dest = $d.Range($paraEnd, $paraEnd)
$dest.FormattedText = $formatted

$newRun = $d.Range($paraEnd, $paraEnd + 1)
$newRun.Text = "."

Write-Output $target.Range.Text
